$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row (42) with the latest mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(42, 1).Value = "Klacht over levering"
$logs.Cells.Item(42, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(42, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Cells.Item(42, 4).Value = "Klacht"
$logs.Cells.Item(42, 6).Value = "2025-06-17 22:02:40"
$logs.Cells.Item(42, 7).Value = "Nee"

# --- Extend the conditional-formatting ranges to cover the new row 42 ---
$catFormats = $logs.Range("D2:D41").FormatConditions
$catFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D42"))

$answeredFormats = $logs.Range("G2:G41").FormatConditions
$answeredFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G42"))

# --- Sheet "Dashboard": bump the "Klacht" tally to reflect the new row ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(6, 2).Value = 3
